$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.645.78"
$ws.Range("E2").Value = "  -2.52%  "
$ws.Range("D3").Value = "1.887.67"
$ws.Range("E3").Value = "  -3.55%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'321.93"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'0.4563"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("E8").Value = "  -3.64%  "
$ws.Range("D9").Value = "'45.59"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'0.07695"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "'0.9591"
$ws.Range("E11").Value = "  -4.28%  "
$ws.Range("D12").Value = "'21.91"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "1.900.30"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("D14").Value = "'6.953"
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").Value = "'5.646"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("D16").Value = "'0.07038"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = "  -0.26%  "
$ws.Range("D18").Value = "'82.72"
$ws.Range("E18").Value = "  -6.59%  "
$ws.Range("D19").Value = "'0.000009477"
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").Value = "'16.65"
$ws.Range("E20").Value = "  -2.84%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "28.616.88"
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "'5.331"
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("D24").Value = "'10.84"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").Value = "2.108.32"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").Value = "'2.061"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'155.15"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'18.93"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("E29").Value = "  -6.51%  "
$ws.Range("D30").Value = "'116.65"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("D31").Value = "'1.813"
$ws.Range("E31").Value = "  -4.16%  "
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("D33").Value = "'0.8451"
$ws.Range("E33").Value = "  -5.77%  "
$ws.Range("D34").Value = "'5.049"
$ws.Range("E34").Value = "  -3.80%  "
$ws.Range("D35").Value = "'1.246"
$ws.Range("E35").Value = "  -6.82%  "
$ws.Range("D36").Value = "'3.053"
$ws.Range("E36").Value = "  -4.29%  "
$ws.Range("D37").Value = "'1.145"
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("D38").Value = "'0.05615"
$ws.Range("E38").Value = "  -3.62%  "
$ws.Range("D39").Value = "'1.002"
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").Value = "'0.02027"
$ws.Range("E40").Value = "  -4.37%  "
$ws.Range("D41").Value = "'7.425"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("D42").Value = "'0.5463"
$ws.Range("E42").Value = "  -5.10%  "
$ws.Range("D43").Value = "'0.1745"
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("D44").Value = "'0.000002908"
$ws.Range("E44").Value = "  -23.08%  "
$ws.Range("D45").Value = "'9.160"
$ws.Range("E45").Value = "  -6.50%  "
$ws.Range("D46").Value = "'2.694"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("E47").Value = "  -4.17%  "
$ws.Range("D48").Value = "'11.21"
$ws.Range("E48").Value = "  -7.80%  "
$ws.Range("D49").Value = "'2.079"
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").Value = "'0.06759"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("D51").Value = "'110.58"
$ws.Range("E51").Value = "  -3.12%  "
